$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "43.565.22"
Set-TextValue "E2" "  +1.18%  "
Set-TextValue "D3" "2.411.79"
Set-TextValue "E3" "  +2.60%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "306.65"
Set-TextValue "E5" "  +1.37%  "
Set-TextValue "D6" "97.19"
Set-TextValue "E6" "  +1.88%  "
Set-TextValue "D7" "0.508"
Set-TextValue "E7" "  +0.54%  "
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.491"
Set-TextValue "E9" "  -1.04%  "
Set-TextValue "D10" "34.97"
Set-TextValue "E10" "  +2.81%  "
Set-TextValue "E11" "  +1.49%  "
Set-TextValue "E12" "  +2.58%  "
Set-TextValue "D13" "18.49"
Set-TextValue "E13" "  -1.08%  "
Set-TextValue "D14" "6.88"
Set-TextValue "E14" "  +2.68%  "
Set-TextValue "D15" "2.783.09"
Set-TextValue "E15" "  +2.29%  "
Set-TextValue "D16" "2.429.65"
Set-TextValue "E16" "  +2.92%  "
Set-TextValue "D17" "0.824"
Set-TextValue "E17" "  +3.89%  "
Set-TextValue "D18" "43.596.21"
Set-TextValue "E18" "  +1.27%  "
Set-TextValue "E19" "  +2.67%  "
Set-TextValue "D20" "12.13"
Set-TextValue "E20" "  -0.38%  "
Set-TextValue "E21" "  +1.74%  "
Set-TextValue "D22" "68.41"
Set-TextValue "E22" "  +0.54%  "
Set-TextValue "D23" "237.69"
Set-TextValue "E23" "  +1.12%  "
Set-TextValue "D24" "2.24"
Set-TextValue "E24" "  +0.89%  "
Set-TextValue "E25" "  +1.07%  "
Set-TextValue "E26" "  +0.05%  "
Set-TextValue "D27" "24.95"
Set-TextValue "E27" "  +1.95%  "
Set-TextValue "E28" "  -0.59%  "
Set-TextValue "D29" "9.44"
Set-TextValue "E29" "  +3.65%  "
Set-TextValue "D30" "32.46"
Set-TextValue "E30" "  +3.81%  "
Set-TextValue "D31" "0.118"
Set-TextValue "E31" "  +16.57%  "
Set-TextValue "D32" "18.44"
Set-TextValue "E32" "  +7.35%  "
Set-TextValue "D33" "5.12"
Set-TextValue "E33" "  +2.06%  "
Set-TextValue "E34" "  +0.09%  "
Set-TextValue "E35" "  +3.90%  "
Set-TextValue "D36" "133.10"
Set-TextValue "E36" "  +28.17%  "
Set-TextValue "E37" "  +3.29%  "
Set-TextValue "E38" "  +6.72%  "
Set-TextValue "E39" "  +0.44%  "
Set-TextValue "E40" "  -1.17%  "
Set-TextValue "D42" "21.32"
Set-TextValue "E42" "  -5.06%  "
Set-TextValue "D43" "1.945.89"
Set-TextValue "E43" "  +0.23%  "
Set-TextValue "E44" "  +1.79%  "
Set-TextValue "D45" "2.17"
Set-TextValue "E45" "  +2.59%  "
Set-TextValue "D46" "2.84"
Set-TextValue "E46" "  +4.08%  "
Set-TextValue "D47" "9.29"
Set-TextValue "E47" "  -1.18%  "
Set-TextValue "D48" "2.636.33"
Set-TextValue "E48" "  +2.04%  "
Set-TextValue "D49" "1.56"
Set-TextValue "E49" "  +4.49%  "
Set-TextValue "D50" "52.65"
Set-TextValue "E50" "  -0.27%  "
Set-TextValue "D51" "72.27"
Set-TextValue "E51" "  +0.25%  "
